$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column is added after the existing "sum" column (G).
# Copy G1's formatting (bold font, border, centered header style) onto H1
# so the new header looks consistent with the other headers.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data row value for the new column.
$ws.Range("H2").Value = 0
